$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: add the missing "Number of affiliations" value
$ws.Range("I7").Value = 0

# Row 8: new collaborator entry (Richard Hugtenburg)
# Written in the same order the author filled the row (email first).
$ws.Range("E8").Value = "r.p.hugtenburg@swansea.ac.uk"
$ws.Range("A8").Value = "Dr."
$ws.Range("B8").Value = "Richard"
$ws.Range("C8").Value = "Hugtenburg"
$ws.Range("D8").Value = "R.P."
$ws.Range("F8").Value = "R.P. Hugtenburg"
$ws.Range("G8").Value = "Swansea-BioMed"
$ws.Range("H8").Value = "Department of Biomedical Sciences, Faculty of Science and Engineering, Swansea University, Singleton Park, Swansea, SA2 8PP, UK"
$ws.Range("I8").Value = 0

# Match the "email address" column styling used in row 7 (Hyperlink style)
$ws.Range("E8").Style = $ws.Range("E7").Style

# Grow Table1 to include the newly added row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:O8"))

# Update the active selection to reflect where editing left off
$ws.Range("I7").Select()
